$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1971.75
$ws.Range("I19").Value = 1975
$ws.Range("J19").Value = 1949
$ws.Range("K19").Value = 1975
$ws.Range("L19").Value = 1949
$ws.Range("M19").Value = -1800
$ws.Range("N19").Value = -2299
$ws.Range("H98").Value = 741.6429000000001
$ws.Range("I98").Value = 419
$ws.Range("J98").Value = 1171.8334
$ws.Range("K98").Value = 419
$ws.Range("L98").Value = 1171.8334
$ws.Range("M98").Value = 1079
$ws.Range("N98").Value = -4167.8334
$ws.Range("H112").Value = 1174.3889
$ws.Range("J112").Value = 1542.5217
$ws.Range("L112").Value = 4627.5651
$ws.Range("N112").Value = -6843.5651
$ws.Range("H113").Value = 2998.5
$ws.Range("I113").Value = 2998.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2998.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 255.5
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 741.6429000000001
$ws.Range("I122").Value = 419
$ws.Range("J122").Value = 1171.8334
$ws.Range("K122").Value = 1257
$ws.Range("L122").Value = 3515.5002
$ws.Range("M122").Value = 1193
$ws.Range("N122").Value = -8415.5002
$ws.Range("H127").Value = 3197.6667
$ws.Range("I127").Value = 796.5
$ws.Range("K127").Value = 2389.5
$ws.Range("M127").Value = 2570.5
$ws.Range("H137").Value = 4224.5454
$ws.Range("I137").Value = 499.5
$ws.Range("J137").Value = 5052.3335
$ws.Range("K137").Value = 1498.5
$ws.Range("L137").Value = 15157.0005
$ws.Range("M137").Value = 1051.5
$ws.Range("N137").Value = -20257.0005
$ws.Range("H138").Value = 2733.8965
$ws.Range("J138").Value = 2953.913
$ws.Range("L138").Value = 8861.739
$ws.Range("N138").Value = -19141.739

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5944.476
$ws.Range("I32").Value = 5944.476
$ws.Range("K32").Value = 5944.476
$ws.Range("M32").Value = -5657.476
$ws.Range("H97").Value = 2208
$ws.Range("I97").Value = 1090
$ws.Range("J97").Value = 4444
$ws.Range("K97").Value = 1090
$ws.Range("L97").Value = 4444
$ws.Range("M97").Value = -594
$ws.Range("N97").Value = -5436

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2035.4117
$ws.Range("I134").Value = 1854.0769
$ws.Range("K134").Value = 5562.2307
$ws.Range("M134").Value = -3027.2307

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2535.75
$ws.Range("I31").Value = 2610.2856
$ws.Range("K31").Value = 2610.2856
$ws.Range("M31").Value = -2315.2856
$ws.Range("H34").Value = 2535.75
$ws.Range("I34").Value = 2610.2856
$ws.Range("K34").Value = 2610.2856
$ws.Range("M34").Value = -2408.2856
$ws.Range("H86").Value = 26145.285
$ws.Range("I86").Value = 11667.667
$ws.Range("J86").Value = 37003.5
$ws.Range("K86").Value = 11667.667
$ws.Range("L86").Value = 37003.5
$ws.Range("M86").Value = -10544.667
$ws.Range("N86").Value = -39249.5
$ws.Range("H89").Value = 26145.285
$ws.Range("I89").Value = 11667.667
$ws.Range("J89").Value = 37003.5
$ws.Range("K89").Value = 58338.335
$ws.Range("L89").Value = 185017.5
$ws.Range("M89").Value = -52722.335
$ws.Range("N89").Value = -196249.5
$ws.Range("H99").Value = 4680
$ws.Range("I99").Value = 5200
$ws.Range("K99").Value = 5200
$ws.Range("M99").Value = -3702
$ws.Range("H109").Value = 61284
$ws.Range("J109").Value = 61284
$ws.Range("L109").Value = 61284
$ws.Range("N109").Value = -63364
$ws.Range("H122").Value = 4109.6665
$ws.Range("I122").Value = 3831.1667
$ws.Range("K122").Value = 11493.5001
$ws.Range("M122").Value = -9043.500100000001
$ws.Range("H126").Value = 4680
$ws.Range("I126").Value = 5200
$ws.Range("K126").Value = 15600
$ws.Range("M126").Value = -13130
$ws.Range("H134").Value = 3674
$ws.Range("I134").Value = 3674
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11022
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -8487
$ws.Range("N134").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 18442.818
$ws.Range("I56").Value = 18442.818
$ws.Range("K56").Value = 18442.818
$ws.Range("M56").Value = -17912.818
$ws.Range("H98").Value = 25759.4
$ws.Range("I98").Value = 6900
$ws.Range("J98").Value = 30474.25
$ws.Range("K98").Value = 20700
$ws.Range("L98").Value = 91422.75
$ws.Range("M98").Value = -19202
$ws.Range("N98").Value = -94418.75
$ws.Range("H132").Value = 5893.615
$ws.Range("I132").Value = 4566.3335
$ws.Range("J132").Value = 6291.8
$ws.Range("K132").Value = 41097.0015
$ws.Range("L132").Value = 56626.2
$ws.Range("M132").Value = -38567.0015
$ws.Range("N132").Value = -61686.2
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H137").Value = 2554.125
$ws.Range("I137").Value = 2480
$ws.Range("J137").Value = 2677.6667
$ws.Range("K137").Value = 7440
$ws.Range("L137").Value = 8033.000100000001
$ws.Range("M137").Value = -2340
$ws.Range("N137").Value = -18233.0001
$ws.Range("H138").Value = 6434.5
$ws.Range("I138").Value = 5246.3335
$ws.Range("K138").Value = 15739.0005
$ws.Range("M138").Value = -10599.0005
$ws.Range("H139").Value = 1530
$ws.Range("I139").Value = 1530
$ws.Range("K139").Value = 4590
$ws.Range("M139").Value = 550
$ws.Range("H140").Value = 3000
$ws.Range("I140").Value = 3000
$ws.Range("K140").Value = 9000
$ws.Range("M140").Value = -3820

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 80.5
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("H107").Value = 844
$ws.Range("I107").Value = 844
$ws.Range("K107").Value = 844
$ws.Range("M107").Value = 1076

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 101277.2
$ws.Range("I62").Value = 125746.5
$ws.Range("J62").Value = 3400
$ws.Range("K62").Value = 125746.5
$ws.Range("L62").Value = 3400
$ws.Range("M62").Value = -125122.5
$ws.Range("N62").Value = -4648
$ws.Range("H65").Value = 101277.2
$ws.Range("I65").Value = 125746.5
$ws.Range("J65").Value = 3400
$ws.Range("K65").Value = 628732.5
$ws.Range("L65").Value = 17000
$ws.Range("M65").Value = -625612.5
$ws.Range("N65").Value = -23240
$ws.Range("H132").Value = 3818.4443
$ws.Range("I132").Value = 3314.1428
$ws.Range("J132").Value = 4139.364
$ws.Range("K132").Value = 9942.428400000001
$ws.Range("L132").Value = 12418.092
$ws.Range("M132").Value = -7412.428400000001
$ws.Range("N132").Value = -17478.092
